$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "id_scenario" column (column A) is no longer needed in the table,
# so delete it entirely -- this shifts every other column one slot to
# the left (B->A, C->B, ... Q->P) and drops the now-unused shared string.
$ws.Range("A1").EntireColumn.Delete()

# Resize the table to match the new (narrower) extent.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:P9"))

# Re-assert the header labels so the table's column definitions stay in
# sync with the shifted header row.
$headerNames = @("id_region","id_sector","id_unit_user_type","unit","2011","2012","2013","2014","2015","2016","2017","2018","2019","2020","2021","2022")
for ($i = 1; $i -le $headerNames.Length; $i++) {
    $ws.Cells.Item(1, $i).Value = $headerNames[$i - 1]
}

# Restore the selection to match the post-edit state.
$ws.Range("C14").Select()
